# Form_Change_Concrete_23_Aug_1pm.xlsx - apply the "23 Aug 1pm" documentation edits:
#   1. Fix a typo in the Purpose cell: "locial" -> "logical"
#   2. Move the sheet's active selection from D3:E3 to A5:E5 (the Purpose row)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Correct the typo in the Purpose row (merged cell A5:E5) -------------
$ws.Range("A5").Value = "Purpose: Unit test the logical structure of the Form_Change_Concrete Class its Interface"

# --- 2. Update the saved selection to the Purpose row (A5:E5) ---------------
$ws.Activate() | Out-Null
$ws.Range("A5:E5").Select() | Out-Null
